# Populate the "Flag_path" column (AK) on the Data sheet with the path to
# each project's flags TSV file, for every IFCB project currently listed.
# (Mirrors running "step2" of the project standardizer, which writes out
# the per-project flag file location.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$flagPaths = @{
    2248 = "~/GIT/PSSdb/raw/flags/IFCB/project_2248_flags.tsv"
    3147 = "~/GIT/PSSdb/raw/flags/IFCB/project_3147_flags.tsv"
    3289 = "~/GIT/PSSdb/raw/flags/IFCB/project_3289_flags.tsv"
    3290 = "~/GIT/PSSdb/raw/flags/IFCB/project_3290_flags.tsv"
    3294 = "~/GIT/PSSdb/raw/flags/IFCB/project_3294_flags.tsv"
    3295 = "~/GIT/PSSdb/raw/flags/IFCB/project_3295_flags.tsv"
    3296 = "~/GIT/PSSdb/raw/flags/IFCB/project_3296_flags.tsv"
    3297 = "~/GIT/PSSdb/raw/flags/IFCB/project_3297_flags.tsv"
    3298 = "~/GIT/PSSdb/raw/flags/IFCB/project_3298_flags.tsv"
    3299 = "~/GIT/PSSdb/raw/flags/IFCB/project_3299_flags.tsv"
    3300 = "~/GIT/PSSdb/raw/flags/IFCB/project_3300_flags.tsv"
    3301 = "~/GIT/PSSdb/raw/flags/IFCB/project_3301_flags.tsv"
    3302 = "~/GIT/PSSdb/raw/flags/IFCB/project_3302_flags.tsv"
    3303 = "~/GIT/PSSdb/raw/flags/IFCB/project_3303_flags.tsv"
    3304 = "~/GIT/PSSdb/raw/flags/IFCB/project_3304_flags.tsv"
    3305 = "~/GIT/PSSdb/raw/flags/IFCB/project_3305_flags.tsv"
    3306 = "~/GIT/PSSdb/raw/flags/IFCB/project_3306_flags.tsv"
    3307 = "~/GIT/PSSdb/raw/flags/IFCB/project_3307_flags.tsv"
    3308 = "~/GIT/PSSdb/raw/flags/IFCB/project_3308_flags.tsv"
    3309 = "~/GIT/PSSdb/raw/flags/IFCB/project_3309_flags.tsv"
    3310 = "~/GIT/PSSdb/raw/flags/IFCB/project_3310_flags.tsv"
    3311 = "~/GIT/PSSdb/raw/flags/IFCB/project_3311_flags.tsv"
    3312 = "~/GIT/PSSdb/raw/flags/IFCB/project_3312_flags.tsv"
    3313 = "~/GIT/PSSdb/raw/flags/IFCB/project_3313_flags.tsv"
    3314 = "~/GIT/PSSdb/raw/flags/IFCB/project_3314_flags.tsv"
    3315 = "~/GIT/PSSdb/raw/flags/IFCB/project_3315_flags.tsv"
    3318 = "~/GIT/PSSdb/raw/flags/IFCB/project_3318_flags.tsv"
    3320 = "~/GIT/PSSdb/raw/flags/IFCB/project_3320_flags.tsv"
    3321 = "~/GIT/PSSdb/raw/flags/IFCB/project_3321_flags.tsv"
    3322 = "~/GIT/PSSdb/raw/flags/IFCB/project_3322_flags.tsv"
    3323 = "~/GIT/PSSdb/raw/flags/IFCB/project_3323_flags.tsv"
    3324 = "~/GIT/PSSdb/raw/flags/IFCB/project_3324_flags.tsv"
    3325 = "~/GIT/PSSdb/raw/flags/IFCB/project_3325_flags.tsv"
    3326 = "~/GIT/PSSdb/raw/flags/IFCB/project_3326_flags.tsv"
    3331 = "~/GIT/PSSdb/raw/flags/IFCB/project_3331_flags.tsv"
    3332 = "~/GIT/PSSdb/raw/flags/IFCB/project_3332_flags.tsv"
    3333 = "~/GIT/PSSdb/raw/flags/IFCB/project_3333_flags.tsv"
    3334 = "~/GIT/PSSdb/raw/flags/IFCB/project_3334_flags.tsv"
    3335 = "~/GIT/PSSdb/raw/flags/IFCB/project_3335_flags.tsv"
    3337 = "~/GIT/PSSdb/raw/flags/IFCB/project_3337_flags.tsv"
    3338 = "~/GIT/PSSdb/raw/flags/IFCB/project_3338_flags.tsv"
    3339 = "~/GIT/PSSdb/raw/flags/IFCB/project_3339_flags.tsv"
    3340 = "~/GIT/PSSdb/raw/flags/IFCB/project_3340_flags.tsv"
    3341 = "~/GIT/PSSdb/raw/flags/IFCB/project_3341_flags.tsv"
    3342 = "~/GIT/PSSdb/raw/flags/IFCB/project_3342_flags.tsv"
    3343 = "~/GIT/PSSdb/raw/flags/IFCB/project_3343_flags.tsv"
}

# Column A holds the Project_ID, column AK (37th column) holds Flag_path.
$firstRow = $ws.UsedRange.Row
$lastRow = $firstRow + $ws.UsedRange.Rows.Count - 1
for ($r = 2; $r -le $lastRow; $r++) {
    $pid = $ws.Cells.Item($r, 1).Value2
    if ($pid -is [double] -and $flagPaths.ContainsKey([int]$pid)) {
        $ws.Cells.Item($r, 37).Value2 = $flagPaths[[int]$pid]
    }
}
